{"js": "const body = context.document.body;\nconst pairs = [\n  [\"2025-04-04 Friday\", \"2025-04-05 Saturday\"],\n  [\"7+15=22\", \"84-29=55\"],\n  [\"37-29=8\", \"25+27=52\"],\n  [\"35+57=92\", \"91-8=83\"],\n  [\"28+36=64\", \"73-44=29\"],\n  [\"27+7=34\", \"71-55=16\"],\n  [\"36+47=83\", \"24+59=83\"],\n  [\"44+19=63\", \"73-19=54\"],\n  [\"25-19=6\", \"25+17=42\"],\n  [\"72-47=25\", \"28+5=33\"],\n  [\"74+7=81\", \"73-24=49\"],\n  [\"6+46=52\", \"28+66=94\"],\n  [\"32+39=71\", \"84-39=45\"],\n  [\"50-24=26\", \"68+18=86\"],\n  [\"32-26=6\", \"41-7=34\"],\n  [\"9+42=51\", \"14+68=82\"],\n  [\"44-25=19\", \"83-6=77\"],\n  [\"51-42=9\", \"12+69=81\"],\n  [\"35+29=64\", \"33-14=19\"],\n  [\"43-15=28\", \"7+55=62\"],\n  [\"37+8=45\", \"71-12=59\"],\n  [\"81-43=38\", \"64-27=37\"],\n  [\"55+16=71\", \"46+39=85\"],\n  [\"85-36=49\", \"50-1=49\"],\n  [\"96-18=78\", \"74-48=26\"],\n  [\"32-13=19\", \"71-5=66\"],\n  [\"63-45=18\", \"8+3=11\"],\n  [\"22+9=31\", \"46+17=63\"],\n  [\"41-19=22\", \"91-68=23\"],\n  [\"24+8=32\", \"76+15=91\"],\n  [\"19+44=63\", \"49+6=55\"],\n  [\"27+35=62\", \"66+19=85\"],\n  [\"68+9=77\", \"55-8=47\"],\n  [\"81-53=28\", \"47+17=64\"],\n  [\"51-7=44\", \"18+74=92\"],\n  [\"72-63=9\", \"46+6=52\"],\n  [\"95-27=68\", \"27+39=66\"],\n  [\"84-45=39\", \"91-46=45\"],\n  [\"23+18=41\", \"67+8=75\"],\n  [\"59+25=84\", \"66+19=85\"],\n  [\"48+7=55\", \"81-79=2\"],\n  [\"84-77=7\", \"46+19=65\"],\n  [\"63-29=34\", \"17+19=36\"],\n  [\"25+58=83\", \"73-59=14\"],\n  [\"88+3=91\", \"18+63=81\"],\n  [\"57+6=63\", \"97-19=78\"],\n  [\"83-46=37\", \"93-25=68\"],\n  [\"25+59=84\", \"39+59=98\"],\n  [\"33+8=41\", \"57+38=95\"],\n  [\"29+22=51\", \"7+84=91\"],\n  [\"82-3=79\", \"5+86=91\"],\n  [\"81-7=74\", \"82-23=59\"],\n  [\"53-8=45\", \"27+67=94\"],\n  [\"85-38=47\", \"91-65=26\"],\n  [\"71-63=8\", \"73-27=46\"],\n  [\"45-8=37\", \"5+38=43\"],\n  [\"5+76=81\", \"60-47=13\"],\n  [\"48+33=81\", \"68+17=85\"],\n  [\"9+27=36\", \"30-14=16\"],\n  [\"20-9=11\", \"34+17=51\"],\n  [\"58+36=94\", \"64-25=39\"],\n  [\"19+65=84\", \"76-19=57\"],\n  [\"46+48=94\", \"17+48=65\"],\n  [\"8+48=56\", \"64-19=45\"],\n  [\"31-29=2\", \"22-8=14\"],\n  [\"82-5=77\", \"60-5=55\"],\n  [\"55+26=81\", \"86+5=91\"],\n  [\"78+17=95\", \"78+6=84\"],\n  [\"6+18=24\", \"83-79=4\"],\n  [\"33-18=15\", \"3+59=62\"],\n  [\"92-67=25\", \"48+49=97\"],\n  [\"42+49=91\", \"19+42=61\"],\n  [\"52-3=49\", \"36+9=45\"],\n  [\"95-19=76\", \"64-48=16\"],\n  [\"35+28=63\", \"41-33=8\"],\n  [\"5+56=61\", \"70-36=34\"],\n  [\"20-13=7\", \"72-36=36\"],\n  [\"66-18=48\", \"18+37=55\"],\n  [\"95-68=27\", \"43-6=37\"],\n  [\"67+6=73\", \"49+7=56\"],\n  [\"66-37=29\", \"29+44=73\"],\n  [\"81-32=49\", \"71-54=17\"],\n  [\"48+16=64\", \"26+29=55\"],\n  [\"29+54=83\", \"90-28=62\"],\n  [\"60-52=8\", \"62-59=3\"],\n  [\"64-49=15\", \"3+78=81\"],\n  [\"72-26=46\", \"44-27=17\"],\n  [\"19+2=21\", \"13+38=51\"],\n  [\"78+19=97\", \"17+16=33\"],\n  [\"94-45=49\", \"91-16=75\"],\n  [\"12+39=51\", \"17+25=42\"],\n  [\"4+18=22\", \"9+17=26\"],\n  [\"71-66=5\", \"58+4=62\"],\n  [\"83-76=7\", \"34+17=51\"],\n  [\"90-41=49\", \"44-7=37\"],\n  [\"12+9=21\", \"90-66=24\"],\n  [\"80-18=62\", \"53-19=34\"],\n  [\"42-23=19\", \"36-8=28\"],\n  [\"52-16=36\", \"78-59=19\"],\n  [\"36+17=53\", \"95-26=69\"],\n  [\"7+5=12\", \"46+5=51\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  ,@(\"2025-04-04 Friday\", \"2025-04-05 Saturday\")\n  ,@(\"7+15=22\", \"84-29=55\")\n  ,@(\"37-29=8\", \"25+27=52\")\n  ,@(\"35+57=92\", \"91-8=83\")\n  ,@(\"28+36=64\", \"73-44=29\")\n  ,@(\"27+7=34\", \"71-55=16\")\n  ,@(\"36+47=83\", \"24+59=83\")\n  ,@(\"44+19=63\", \"73-19=54\")\n  ,@(\"25-19=6\", \"25+17=42\")\n  ,@(\"72-47=25\", \"28+5=33\")\n  ,@(\"74+7=81\", \"73-24=49\")\n  ,@(\"6+46=52\", \"28+66=94\")\n  ,@(\"32+39=71\", \"84-39=45\")\n  ,@(\"50-24=26\", \"68+18=86\")\n  ,@(\"32-26=6\", \"41-7=34\")\n  ,@(\"9+42=51\", \"14+68=82\")\n  ,@(\"44-25=19\", \"83-6=77\")\n  ,@(\"51-42=9\", \"12+69=81\")\n  ,@(\"35+29=64\", \"33-14=19\")\n  ,@(\"43-15=28\", \"7+55=62\")\n  ,@(\"37+8=45\", \"71-12=59\")\n  ,@(\"81-43=38\", \"64-27=37\")\n  ,@(\"55+16=71\", \"46+39=85\")\n  ,@(\"85-36=49\", \"50-1=49\")\n  ,@(\"96-18=78\", \"74-48=26\")\n  ,@(\"32-13=19\", \"71-5=66\")\n  ,@(\"63-45=18\", \"8+3=11\")\n  ,@(\"22+9=31\", \"46+17=63\")\n  ,@(\"41-19=22\", \"91-68=23\")\n  ,@(\"24+8=32\", \"76+15=91\")\n  ,@(\"19+44=63\", \"49+6=55\")\n  ,@(\"27+35=62\", \"66+19=85\")\n  ,@(\"68+9=77\", \"55-8=47\")\n  ,@(\"81-53=28\", \"47+17=64\")\n  ,@(\"51-7=44\", \"18+74=92\")\n  ,@(\"72-63=9\", \"46+6=52\")\n  ,@(\"95-27=68\", \"27+39=66\")\n  ,@(\"84-45=39\", \"91-46=45\")\n  ,@(\"23+18=41\", \"67+8=75\")\n  ,@(\"59+25=84\", \"66+19=85\")\n  ,@(\"48+7=55\", \"81-79=2\")\n  ,@(\"84-77=7\", \"46+19=65\")\n  ,@(\"63-29=34\", \"17+19=36\")\n  ,@(\"25+58=83\", \"73-59=14\")\n  ,@(\"88+3=91\", \"18+63=81\")\n  ,@(\"57+6=63\", \"97-19=78\")\n  ,@(\"83-46=37\", \"93-25=68\")\n  ,@(\"25+59=84\", \"39+59=98\")\n  ,@(\"33+8=41\", \"57+38=95\")\n  ,@(\"29+22=51\", \"7+84=91\")\n  ,@(\"82-3=79\", \"5+86=91\")\n  ,@(\"81-7=74\", \"82-23=59\")\n  ,@(\"53-8=45\", \"27+67=94\")\n  ,@(\"85-38=47\", \"91-65=26\")\n  ,@(\"71-63=8\", \"73-27=46\")\n  ,@(\"45-8=37\", \"5+38=43\")\n  ,@(\"5+76=81\", \"60-47=13\")\n  ,@(\"48+33=81\", \"68+17=85\")\n  ,@(\"9+27=36\", \"30-14=16\")\n  ,@(\"20-9=11\", \"34+17=51\")\n  ,@(\"58+36=94\", \"64-25=39\")\n  ,@(\"19+65=84\", \"76-19=57\")\n  ,@(\"46+48=94\", \"17+48=65\")\n  ,@(\"8+48=56\", \"64-19=45\")\n  ,@(\"31-29=2\", \"22-8=14\")\n  ,@(\"82-5=77\", \"60-5=55\")\n  ,@(\"55+26=81\", \"86+5=91\")\n  ,@(\"78+17=95\", \"78+6=84\")\n  ,@(\"6+18=24\", \"83-79=4\")\n  ,@(\"33-18=15\", \"3+59=62\")\n  ,@(\"92-67=25\", \"48+49=97\")\n  ,@(\"42+49=91\", \"19+42=61\")\n  ,@(\"52-3=49\", \"36+9=45\")\n  ,@(\"95-19=76\", \"64-48=16\")\n  ,@(\"35+28=63\", \"41-33=8\")\n  ,@(\"5+56=61\", \"70-36=34\")\n  ,@(\"20-13=7\", \"72-36=36\")\n  ,@(\"66-18=48\", \"18+37=55\")\n  ,@(\"95-68=27\", \"43-6=37\")\n  ,@(\"67+6=73\", \"49+7=56\")\n  ,@(\"66-37=29\", \"29+44=73\")\n  ,@(\"81-32=49\", \"71-54=17\")\n  ,@(\"48+16=64\", \"26+29=55\")\n  ,@(\"29+54=83\", \"90-28=62\")\n  ,@(\"60-52=8\", \"62-59=3\")\n  ,@(\"64-49=15\", \"3+78=81\")\n  ,@(\"72-26=46\", \"44-27=17\")\n  ,@(\"19+2=21\", \"13+38=51\")\n  ,@(\"78+19=97\", \"17+16=33\")\n  ,@(\"94-45=49\", \"91-16=75\")\n  ,@(\"12+39=51\", \"17+25=42\")\n  ,@(\"4+18=22\", \"9+17=26\")\n  ,@(\"71-66=5\", \"58+4=62\")\n  ,@(\"83-76=7\", \"34+17=51\")\n  ,@(\"90-41=49\", \"44-7=37\")\n  ,@(\"12+9=21\", \"90-66=24\")\n  ,@(\"80-18=62\", \"53-19=34\")\n  ,@(\"42-23=19\", \"36-8=28\")\n  ,@(\"52-16=36\", \"78-59=19\")\n  ,@(\"36+17=53\", \"95-26=69\")\n  ,@(\"7+5=12\", \"46+5=51\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
